$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing table body (rows 5-21) before rewriting with final layout
$ws.Range("B5:H21").ClearContents()

$data = New-Object 'object[,]' 18,7
$data[0,0] = 'The Combinatorial Theory of Single-Elimination Tournaments'
$data[0,1] = 'Christopher Todd Edwards'
$data[0,2] = 1991
$data[0,3] = 'https://scholarworks.montana.edu/xmlui/bitstream/handle/1/6870/31762100987518.pdf?sequence=1'
$data[0,4] = 'Ordered Brackets'
$data[0,5] = 'x'
$data[0,6] = 'Edward''s Theorem'

$data[1,0] = 'The Method of Paired Comparisons'
$data[1,1] = 'H.A. David'
$data[1,2] = 1963
$data[1,4] = 'Ordered Brackets'
$data[1,5] = 'x'
$data[1,6] = 'SST'

$data[2,0] = 'When 15th Is Better Than 8th: The Math Shows the Bracket Is Backward'
$data[2,1] = 'Nate Silver'
$data[2,2] = 2011
$data[2,3] = 'https://fivethirtyeight.com/features/when-15th-is-better-than-8th-the-math-shows-the-bracket-is-backward/'
$data[2,4] = 'Ordered Brackets'
$data[2,5] = 'x'
$data[2,6] = 'March Madness not ordered'

$data[3,0] = 'Comparing Draws for Single Elimination Tournaments'
$data[3,1] = 'Jeff Horen and Raymond Riezman'
$data[3,2] = 1984
$data[3,3] = 'https://sci-hub.se/https://www.jstor.org/stable/170742'
$data[3,4] = 'Ordered Brackets'
$data[3,5] = 'x'
$data[3,6] = 'Defining Ordered (named fairness), proving 8-balanced is not ordered'

$data[4,0] = 'Do Stronger Players Win More Knockout Tournaments'
$data[4,1] = 'Robert Chen and F. K. Hwang'
$data[4,2] = 1978
$data[4,3] = 'https://www-jstor-org.ezp-prod1.hul.harvard.edu/stable/pdf/2286606.pdf?refreqid=excelsior%3Ac6f9f27cd7701f5e7471bf45a00f7c81&ab_segments=&origin=&initiator=&acceptTC=1'
$data[4,4] = 'Ordered Brackets'
$data[4,5] = 'x'
$data[4,6] = 'defined ordered'

$data[5,0] = 'The UPA Manual of Championship Series Tournament Formats'
$data[5,1] = 'Eric Simon'
$data[5,2] = 2008
$data[5,3] = 'https://usaultimate.org/wp-content/uploads/2020/11/USAU_TournamentFormats.pdf'
$data[5,4] = 'Semibrackets'
$data[5,5] = 'x'
$data[5,6] = 'ultimate formats'

$data[6,0] = 'New Concepts in Seeding Knockout Tournaments'
$data[6,1] = 'F. K. Hwang'
$data[6,2] = 1982
$data[6,3] = 'https://www-jstor-org.ezp-prod1.hul.harvard.edu/stable/pdf/2320220.pdf?refreqid=excelsior%3A5ecd05b3217d589d8ca6170e0208a599&ab_segments=&origin=&initiator=&acceptTC=1'
$data[6,4] = 'Tiered Seedings'
$data[6,5] = 'x'
$data[6,6] = 'reseeding is ordered'

$data[7,0] = 'Efficacy of traditional sport tournament structures'
$data[7,1] = 'T McGarry and RW Schutz'
$data[7,2] = 1997
$data[7,3] = 'https://www-jstor-org.ezp-prod1.hul.harvard.edu/stable/3009945?sid=primo&origin=crossref'
$data[7,4] = 'Swiss Systems'
$data[7,5] = 'x'
$data[7,6] = 'Swapping 4th and 5th in swiss vs in 3rd place games'

$data[8,0] = 'Stronger Players Win More Balanced Knockout Tournaments'
$data[8,1] = 'Robert Chen and F. K. Hwang'
$data[8,2] = ' '
$data[8,3] = 'https://link-springer-com.ezp-prod1.hul.harvard.edu/content/pdf/10.1007/BF01864157.pdf'
$data[8,4] = 'Tiered Seedings'
$data[8,5] = '?'
$data[8,6] = 'fully randomized balanced brackets are ordered'

$data[9,0] = 'Stronger Players Need not Win More Knockout Tournaments'
$data[9,1] = 'Robert Israel'
$data[9,2] = 1981
$data[9,3] = 'https://www-jstor-org.ezp-prod1.hul.harvard.edu/stable/2287594?sid=primo'
$data[9,4] = 'Tiered Seedings'
$data[9,5] = '?'
$data[9,6] = 'example of fully randomized bracket that is not ordered'

$data[10,0] = 'What is the Correct Way to Seed a Knockout Tournament'
$data[10,1] = 'Allen Schwenk'
$data[10,2] = 2018
$data[10,3] = 'https://www-tandfonline-com.ezp-prod1.hul.harvard.edu/doi/abs/10.1080/00029890.2000.12005171'
$data[10,4] = 'Tiered Seedings'
$data[10,5] = '?'
$data[10,6] = 'cohort randomized'

$data[11,0] = 'OEIS Sequence A002572'
$data[11,1] = 'OEIS Foundation'
$data[11,2] = 2012
$data[11,3] = 'https://oeis.org/A002572'
$data[11,4] = 'Bracket Signatures'
$data[11,6] = 'Partitions of n'

$data[12,0] = 'Fair Seeding in Knockout Tournametns'
$data[12,1] = 'Thuc Vu and Yoav Shoham'
$data[12,2] = 2011
$data[12,3] = 'https://dl-acm-org.ezp-prod1.hul.harvard.edu/doi/pdf/10.1145%2F2036264.2036273'
$data[12,4] = 'Ordered Brackets'
$data[12,6] = 'Bad version of ordered brackets paper'

$data[13,0] = 'Designing fair 8- and 16-team knockout tournaments'
$data[13,1] = 'Mike Price, J. Cole Prince and Jospeh Geunes'
$data[13,2] = 2011
$data[13,3] = 'https://watermark-silverchair-com.ezp-prod1.hul.harvard.edu/dpr024.pdf?token=AQECAHi208BE49Ooan9kkhW_Ercy7Dm3ZL_9Cf3qfKAc485ysgAAA2MwggNfBgkqhkiG9w0BBwagggNQMIIDTAIBADCCA0UGCSqGSIb3DQEHATAeBglghkgBZQMEAS4wEQQMM48srvUHO6hh9hDhAgEQgIIDFnCTAZLVDcfxa42odq-yePFs3PspALWxR2DPtTlcYSqW_4k0Az8Fzlu-0UizO2rHjtX_VdUhj2B6qZu-C23Jfch07VBUV9LdenqEhWkZKNeWc-a2CvLPOKgi4k43RoV7ATMr7v-z4oBgYvIqTaUBF-j4mUj6ecMouPouJQNYRiy40WukfALvbKGjMudcONZFg1XLRuKuo8z5dV7-hn3i2PedTdn-Aqd97xsJnBKU2rSKpmu1J_B71Zc30ZDXdXOOgTEVCDObA9WiCkBoMu0tk_S5wV8DbfLOANZHnsBsIttOqfAsCIzJR4I4udHER7SV1mYTxuebxzs5g3P8PaWFoSe6y0L8ZfYnqP-2XFc-33v-G71RD7RrXX0x2BGur-6MLA1tXUQrScsMtZivpIxpQNU-3eaIkyDImT_p0kDSaGbO6pJoMY8Ry5ootWCTN9567zYnNAxE8VxWXe9cPSXXPqZSFZl4w9A74qREWKcWb-ktvjSbHv7bQ7BBgNFwZ4bSsa5RQWmMnYcHVflSdYRVa8HdEKoU6ZUGoYV1mvJ-qWA5Q3flRmInc32DqP7ka5UYUp6uiOMfGy7SOP6WwUs_wNsqJe67PnbMpJ214r3sjPsXFIpe6ZD121lIVews1wFeLIOXON1naNFLqzakYkmxtudoLSXiJXmP4eVNVBNQxwU1C3Ig3WjPoIIwqx_88JcCsC-Zxhn749wnNSFTbri_Xj--BorGTcFivud4WvYzTU2JO4ui65Ij9DFQIylIEAGCp16XVdIQsddevXj7Eim_pQ-6VwAZ3KyU5tbLvVCwYF5JgglhVxyUNjCqLoSkUXI8XKX1-lIgWlxY6vS1nchgr-RIM1atFHPrTlT8OXtNlh1DC8tal5sr_67aCTLNiUFqk2Y9NF2BkBhnKwlEIihIxiS9s5PNl5rkHIlQxXozXc_YNmEreoSSFXDXnMQ1u-dZRFUtDXwb0SUhVjXmlTDkLA9Vwn5gtlX-QlId7J82NVqHBV1j22wIuCIJDsVpzJTfDUi1rEy0lKXvASIprCfh4gAblTBCdPQ'
$data[13,4] = 'Ordered Brackets'
$data[13,6] = 'Analysis of 8 team brackets (23 possible configs)'

$data[14,0] = 'Bayesian optimal design of fixed knockout'
$data[14,1] = 'Jonathan Hennessy EMAIL logo and Mark Glickman'
$data[14,2] = 2015
$data[14,4] = 'Simulation'
$data[14,6] = 'simulated annealing of brackets'

$data[15,0] = 'Efficient Simulation of a Random Knockout Tournament'
$data[15,1] = 'Sheldon M. Ross, Samim Ghamami'
$data[15,2] = 2008
$data[15,3] = 'https://www.jise.ir/article_3970_02d5bdcdb9c2ffe9b8b030153c8cf883.pdf'
$data[15,4] = 'Simulation'

$data[16,0] = 'Random Knockout Tournaments'
$data[16,1] = 'Ilan Adler, Yang Cao, Richard Karp, Erol A. Peköz, Sheldon M. Ross'
$data[16,2] = 2017
$data[16,3] = 'https://pubsonline-informs-org.ezp-prod1.hul.harvard.edu/doi/pdf/10.1287%2Fopre.2017.1657'
$data[16,4] = 'Tiered Seedings'
$data[16,6] = 'randomize over bracket shapes'

$data[17,0] = 'Anomalies in Tournament Design: The Madness of March Madness'
$data[17,1] = 'Robert Baumann, Victor Matheson, Cara Howe'
$data[17,2] = 2010
$data[17,3] = 'https://sci-hub.se/https://doi.org/10.2202/1559-0410.1233'
$data[17,4] = 'Tiered Seedings'
$data[17,5] = '?'
$data[17,6] = 'downisdes of reseeding'

$ws.Range("B5:H22").Value2 = $data

[void]$ws.Range("I25").Select()